$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every value as text (inline/shared strings), even
# though the strings look numeric. For each edited row we briefly switch the
# C:D cells to Text format so the new value is stored as text (not coerced to
# a number), write the value, then clear the temporary format again so the
# cell keeps the workbook default style (matching every other data cell).

$row2 = $ws.Range("C2:D2")
$row2.NumberFormat = "@"
$ws.Range("C2").Value = "432"
$ws.Range("D2").Value = "1009828.79"
$row2.ClearFormats()

$row3 = $ws.Range("C3:D3")
$row3.NumberFormat = "@"
$ws.Range("C3").Value = "6"
$ws.Range("D3").Value = "26216.00"
$row3.ClearFormats()

$row4 = $ws.Range("C4:D4")
$row4.NumberFormat = "@"
$ws.Range("C4").Value = "891"
$ws.Range("D4").Value = "2898127.87"
$row4.ClearFormats()

$row6 = $ws.Range("C6:D6")
$row6.NumberFormat = "@"
$ws.Range("C6").Value = "554"
$ws.Range("D6").Value = "1633600.55"
$row6.ClearFormats()

$row7 = $ws.Range("C7:D7")
$row7.NumberFormat = "@"
$ws.Range("C7").Value = "13"
$ws.Range("D7").Value = "29000.00"
$row7.ClearFormats()

$row8 = $ws.Range("C8:D8")
$row8.NumberFormat = "@"
$ws.Range("C8").Value = "26"
$ws.Range("D8").Value = "94093.58"
$row8.ClearFormats()

$row14 = $ws.Range("C14:D14")
$row14.NumberFormat = "@"
$ws.Range("C14").Value = "211"
$ws.Range("D14").Value = "564362.00"
$row14.ClearFormats()

$row16 = $ws.Range("C16:D16")
$row16.NumberFormat = "@"
$ws.Range("C16").Value = "473"
$ws.Range("D16").Value = "1706154.75"
$row16.ClearFormats()

$row19 = $ws.Range("C19:D19")
$row19.NumberFormat = "@"
$ws.Range("C19").Value = "7"
$ws.Range("D19").Value = "19876.16"
$row19.ClearFormats()

$row20 = $ws.Range("C20:D20")
$row20.NumberFormat = "@"
$ws.Range("C20").Value = "167"
$ws.Range("D20").Value = "424599.00"
$row20.ClearFormats()

$row21 = $ws.Range("C21:D21")
$row21.NumberFormat = "@"
$ws.Range("C21").Value = "326"
$ws.Range("D21").Value = "1139929.00"
$row21.ClearFormats()

$row24 = $ws.Range("C24:D24")
$row24.NumberFormat = "@"
$ws.Range("C24").Value = "9"
$ws.Range("D24").Value = "42700.00"
$row24.ClearFormats()

$row28 = $ws.Range("C28:D28")
$row28.NumberFormat = "@"
$ws.Range("C28").Value = "260"
$ws.Range("D28").Value = "660542.64"
$row28.ClearFormats()

$row30 = $ws.Range("C30:D30")
$row30.NumberFormat = "@"
$ws.Range("C30").Value = "517"
$ws.Range("D30").Value = "2051800.70"
$row30.ClearFormats()

$row32 = $ws.Range("C32:D32")
$row32.NumberFormat = "@"
$ws.Range("C32").Value = "367"
$ws.Range("D32").Value = "1211198.17"
$row32.ClearFormats()

$row40 = $ws.Range("C40:D40")
$row40.NumberFormat = "@"
$ws.Range("C40").Value = "130"
$ws.Range("D40").Value = "362102.22"
$row40.ClearFormats()

$row41 = $ws.Range("C41:D41")
$row41.NumberFormat = "@"
$ws.Range("C41").Value = "77"
$ws.Range("D41").Value = "370909.98"
$row41.ClearFormats()

$row42 = $ws.Range("C42:D42")
$row42.NumberFormat = "@"
$ws.Range("C42").Value = "117"
$ws.Range("D42").Value = "474188.99"
$row42.ClearFormats()

$row45 = $ws.Range("C45:D45")
$row45.NumberFormat = "@"
$ws.Range("C45").Value = "342"
$ws.Range("D45").Value = "938867.74"
$row45.ClearFormats()

$row51 = $ws.Range("C51:D51")
$row51.NumberFormat = "@"
$ws.Range("C51").Value = "3419"
$ws.Range("D51").Value = "7815692.31"
$row51.ClearFormats()

$row53 = $ws.Range("C53:D53")
$row53.NumberFormat = "@"
$ws.Range("C53").Value = "3831"
$ws.Range("D53").Value = "13020069.22"
$row53.ClearFormats()

$row55 = $ws.Range("C55:D55")
$row55.NumberFormat = "@"
$ws.Range("C55").Value = "3912"
$ws.Range("D55").Value = "12006253.47"
$row55.ClearFormats()

$row57 = $ws.Range("C57:D57")
$row57.NumberFormat = "@"
$ws.Range("C57").Value = "82"
$ws.Range("D57").Value = "298436.47"
$row57.ClearFormats()

$row73 = $ws.Range("C73:D73")
$row73.NumberFormat = "@"
$ws.Range("C73").Value = "375"
$ws.Range("D73").Value = "926635.70"
$row73.ClearFormats()

$row74 = $ws.Range("C74:D74")
$row74.NumberFormat = "@"
$ws.Range("C74").Value = "5"
$ws.Range("D74").Value = "22254.00"
$row74.ClearFormats()

$row75 = $ws.Range("C75:D75")
$row75.NumberFormat = "@"
$ws.Range("C75").Value = "887"
$ws.Range("D75").Value = "2974751.39"
$row75.ClearFormats()

$row76 = $ws.Range("C76:D76")
$row76.NumberFormat = "@"
$ws.Range("C76").Value = "503"
$ws.Range("D76").Value = "1628543.87"
$row76.ClearFormats()

$row77 = $ws.Range("C77:D77")
$row77.NumberFormat = "@"
$ws.Range("C77").Value = "36"
$ws.Range("D77").Value = "96000.00"
$row77.ClearFormats()
